$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper XML package wrapper used with Range.InsertXML
# ---------------------------------------------------------------------------
function New-PackageXml([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
        '<w:body>' + $bodyFragment + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Title paragraph: wrap the second "linux" run (in "novaordis-linux")
#    with spellStart / spellEnd proofErr markers.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Range = $d.Range($p1.Range.Start, $p1.Range.End)

$p1Fragment = '<w:p w14:paraId="08288892" w14:textId="7041861F" w:rsidR="002D4287" w:rsidRPr="001C62A4" w:rsidRDefault="00544B82" w:rsidP="001D2644">' +
    '<w:pPr><w:pStyle w:val="Title"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve">NovaOrdis </w:t></w:r>' +
    '<w:r w:rsidR="001C62A4" w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>Linux</w:t></w:r>' +
    '<w:r w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r w:rsidR="00ED7747" w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>novaordis-</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidR="001C62A4" w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>linux</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>)</w:t></w:r>' +
    '<w:r w:rsidR="001D2644" w:rsidRPr="001C62A4"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve"> TODO</w:t></w:r>' +
    '</w:p>'

$p1Range.InsertXML((New-PackageXml $p1Fragment))

# ---------------------------------------------------------------------------
# 2) Remove the "_GoBack" bookmark that currently sits in the hyperlink
#    paragraph (it will be re-created at the end of the document below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Last (list) paragraph: add the new "Format instances thread safety..."
#    text together with its proofErr/gramErr markers, and append the
#    "_GoBack" bookmark at the very end of the paragraph / document.
#
#    NOTE: the engine mishandles InsertXML replacements that target the
#    very last paragraph of the document (it appends a stray sibling
#    paragraph instead of replacing in place). Work around this by first
#    adding a temporary empty paragraph after it (so our target paragraph
#    is no longer document-final), performing the replacement, and then
#    removing the temporary trailing paragraph mark again.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$pLast.Range.InsertParagraphAfter()

$pLast = $d.Paragraphs.Item($lastIndex)
$pLastRange = $d.Range($pLast.Range.Start, $pLast.Range.End)

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
$rPrBold = '<w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

$pLastFragment = '<w:p w14:paraId="4E33AE33" w14:textId="77777777" w:rsidR="001D2644" w:rsidRPr="001C62A4" w:rsidRDefault="001D2644" w:rsidP="001D2644">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr>' +
    '<w:r>' + $rPrBold + '<w:t>Format instances thread safety</w:t></w:r>' +
    '<w:r>' + $rPrPlain + '<w:t xml:space="preserve">. Remove public static final </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPrPlain + '<w:t>SimpleDateFormat</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> and other Formats everywhere, and replace them with static method that create the instances on demand. We need to do this because Format instances are not thread safe. See </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPrPlain + '<w:t>io.novaordis.events.api.event</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPrPlain + '<w:t>.DateProperty</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrPlain + '<w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPrPlain + '<w:t>getDefaultDateFormat</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrPlain + '<w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPrPlain + '<w:t>) as example.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$pLastRange.InsertXML((New-PackageXml $pLastFragment))

# Remove the temporary trailing empty paragraph added above.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing.Range.Delete()

Write-Output "edit complete"
